# Fix issue #51 move keep variables after meta variables.
#
# The "summary" worksheet holds, for each enumerator row, pairs of columns
# where the first column of the pair was being populated with the "keep"
# value and the second with the "meta" value (or vice versa). The fix is to
# swap the two values in each pair so the meta variable lines up after the
# keep variable, for every data row (rows 4 through 40) and for every
# column-pair: B/C, D/E, F/G, H/I, J/K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

$firstRow = 4
$lastRow = 40
$pairs = @(("B","C"), ("D","E"), ("F","G"), ("H","I"), ("J","K"))

for ($r = $firstRow; $r -le $lastRow; $r++) {
    foreach ($pair in $pairs) {
        $leftAddr = "$($pair[0])$r"
        $rightAddr = "$($pair[1])$r"

        $leftCell = $ws.Range($leftAddr)
        $rightCell = $ws.Range($rightAddr)

        $leftValue = $leftCell.Value()
        $rightValue = $rightCell.Value()

        $leftCell.Value = $rightValue
        $rightCell.Value = $leftValue
    }
}
